$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header additions: new antioxidant-defense genes Ucp2 (H) and Ucp3 (I) ---
$ws.Range("H1").Value = "Ucp2"
$ws.Range("I1").Value = "Ucp3"

# --- Fill H/I columns for the existing rows 2-9 (row numbers unaffected, still above the new row 10) ---
$ws.Range("H2").Value = 1.443986051446283
$ws.Range("I2").Value = 0.12003842953401894
$ws.Range("H3").Value = 0.82935220096750395
$ws.Range("I3").Value = 0.20899904489604426
$ws.Range("H4").Value = 3.2830957048003695
$ws.Range("I4").Value = 0.40516129193378075
$ws.Range("H5").Value = 1.6472468937100353
$ws.Range("H6").Value = 0.81510442655970228
$ws.Range("I6").Value = 0.40236263846101239
$ws.Range("H7").Value = 3.7452364557597302
$ws.Range("I7").Value = 0.39408213509934548
$ws.Range("H8").Value = 3.6808955351432706
$ws.Range("I8").Value = 1.5226330034909059
$ws.Range("H9").Value = 5.2417800067931557
$ws.Range("I9").Value = 0.91800139715481377

# --- Insert a brand-new data row at position 10 (mouse 2500); shifts old rows 10-22 down to 11-23 ---
$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = 2500
$ws.Range("B10").Value = "MCP230"
$ws.Range("H10").Value = 0.83512080171598624

# --- Fill C-I for the shifted rows 11-23 (these rows already carried A/B/C-G data pre-insert; only H/I are new) ---
$ws.Range("H11").Value = 0.86158008381504336
$ws.Range("H12").Value = 1.1096131622969339
$ws.Range("H13").Value = 2.5492216284270843
$ws.Range("I13").Value = 0.63137465544857885
$ws.Range("H14").Value = 0.42192451104255124
$ws.Range("I14").Value = 0.89289860067401827
$ws.Range("H15").Value = 0.72449960393406643
$ws.Range("I15").Value = 0.33834530031389087
$ws.Range("H16").Value = 0.99313159928742256
$ws.Range("I16").Value = 0.85356223640033568
$ws.Range("H17").Value = 0.2124296140496324
$ws.Range("I17").Value = 0.69090930859305077
$ws.Range("H18").Value = 0.1731455009948227
$ws.Range("I18").Value = 1.3580798004238133
$ws.Range("H19").Value = 0.30567210504475412
$ws.Range("I19").Value = 0.35393792310794098
$ws.Range("H20").Value = 0.036907473875207417
$ws.Range("I20").Value = 0.3859720420326932
$ws.Range("H21").Value = 0.058926353678156636
$ws.Range("I21").Value = 1.6662077050518229
$ws.Range("H22").Value = 1.068107464019588
$ws.Range("I22").Value = 2.4564387291711571
$ws.Range("H23").Value = 4.4848408995347366
$ws.Range("I23").Value = 1.3722736987826982

# --- Leave the selection where the author left it when saving ---
$ws.Range("K18").Select()
